$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Roadmap")
$ws.Activate()

# --- Row 11: G11 status flips from "进行中" to "已完成" ---
$ws.Range("G11").Value = "已完成"

# --- Row 12: becomes an "in progress -> completed" row (style flips from
#     the "in progress" fill to the "completed" fill), gains a completion
#     date in column F, and its status in G flips to "已完成". ---
$ws.Range("A9:G9").Copy()
$ws.Range("A12:G12").PasteSpecial(-4122)
$ws.Range("F12").Value = "2023.11.22"
$ws.Range("G12").Value = "已完成"

# --- New row 14: a new "Material新增默认标准材质" task, styled the same
#     as the other "in progress" rows (e.g. row 8). ---
$ws.Range("A8:G8").Copy()
$ws.Range("A14:G14").PasteSpecial(-4122)
$ws.Range("A14").Value = "Sean Duan"
$ws.Range("B14").Value = "Renderer"
$ws.Range("C14").Value = "Material新增默认标准材质"
$ws.Range("D14").ClearContents()
$ws.Range("E14").Value = "2023.11.22"
$ws.Range("F14").ClearContents()
$ws.Range("G14").Value = "进行中"

# --- Selection moved to C19 (matches the saved sheetView state) ---
$ws.Range("C19").Select()
